$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037579597792125
$ws.Range("D2").Value = 1.040533300745849
$ws.Range("E2").Value = 1.045203871433052
$ws.Range("F2").Value = 1.053372110736914
$ws.Range("I2").Value = 1.038082771872998
$ws.Range("J2").Value = 1.042681699862227
$ws.Range("K2").Value = 1.043315379284932
$ws.Range("L2").Value = 1.047972769223768
$ws.Range("M2").Value = 1.056118261612939
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038466054398504
$ws.Range("D3").Value = 1.041199342010531
$ws.Range("E3").Value = 1.046026717697068
$ws.Range("F3").Value = 1.054356577663086
$ws.Range("I3").Value = 1.038279962851485
$ws.Range("J3").Value = 1.043212910844094
$ws.Range("K3").Value = 1.043792243320464
$ws.Range("L3").Value = 1.048606978714172
$ws.Range("M3").Value = 1.056915320911727
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039040177727296
$ws.Range("D4").Value = 1.041630689919628
$ws.Range("E4").Value = 1.046560055103997
$ws.Range("F4").Value = 1.054994846577108
$ws.Range("I4").Value = 1.038406602487467
$ws.Range("J4").Value = 1.043556519681174
$ws.Range("K4").Value = 1.044100498083387
$ws.Range("L4").Value = 1.049017604598287
$ws.Range("M4").Value = 1.05743170564266
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039281663937199
$ws.Range("D5").Value = 1.041812116719799
$ws.Range("E5").Value = 1.046784484269654
$ws.Range("F5").Value = 1.055263473244371
$ws.Range("I5").Value = 1.038459612283472
$ws.Range("J5").Value = 1.043700942959237
$ws.Range("K5").Value = 1.044230013493025
$ws.Range("L5").Value = 1.049190290323007
$ws.Range("M5").Value = 1.057648944199867
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039322217812475
$ws.Range("D6").Value = 1.041842584201065
$ws.Range("E6").Value = 1.046822179425273
$ws.Range("F6").Value = 1.055308594298051
$ws.Range("I6").Value = 1.03846849938534
$ws.Range("J6").Value = 1.043725190479064
$ws.Range("K6").Value = 1.044251755275923
$ws.Range("L6").Value = 1.04921928841677
$ws.Range("M6").Value = 1.057685428289927
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039043403989255
$ws.Range("D7").Value = 1.041633113809041
$ws.Range("E7").Value = 1.046563053098084
$ws.Range("F7").Value = 1.05499843480915
$ws.Range("I7").Value = 1.03840731170952
$ws.Range("J7").Value = 1.04355844959086
$ws.Range("K7").Value = 1.044102228970307
$ws.Range("L7").Value = 1.049019911804663
$ws.Range("M7").Value = 1.057434607805585
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037879070119658
$ws.Range("D8").Value = 1.040758313994126
$ws.Range("E8").Value = 1.045481768723317
$ws.Range("F8").Value = 1.053704555571891
$ws.Range("I8").Value = 1.03814961108862
$ws.Range("J8").Value = 1.042861249322774
$ws.Range("K8").Value = 1.043476600850453
$ws.Range("L8").Value = 1.048187050798431
$ws.Range("M8").Value = 1.056387499416241
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035831456565464
$ws.Range("D9").Value = 1.039219743367257
$ws.Range("E9").Value = 1.043583375320005
$ws.Range("F9").Value = 1.051434240041506
$ws.Range("I9").Value = 1.037688214077042
$ws.Range("J9").Value = 1.04163182166446
$ws.Range("K9").Value = 1.042371853577857
$ws.Range("L9").Value = 1.046721411990372
$ws.Range("M9").Value = 1.05454727793117
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034469209463584
$ws.Range("D10").Value = 1.038196097926452
$ws.Range("E10").Value = 1.042322556366991
$ws.Range("F10").Value = 1.049927286629879
$ws.Range("I10").Value = 1.037375747318874
$ws.Range("J10").Value = 1.040811682633901
$ws.Range("K10").Value = 1.041633868408017
$ws.Range("L10").Value = 1.045745716251056
$ws.Range("M10").Value = 1.053323851674419
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033880028501522
$ws.Range("D11").Value = 1.037753358224991
$ws.Range("E11").Value = 1.04177776008618
$ws.Range("F11").Value = 1.049276341558851
$ws.Range("I11").Value = 1.03723929857591
$ws.Range("J11").Value = 1.040456443611392
$ws.Range("K11").Value = 1.041313973936863
$ws.Range("L11").Value = 1.04532357585312
$ws.Range("M11").Value = 1.052794915357921
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033661284047918
$ws.Range("D12").Value = 1.037588982503987
$ws.Range("E12").Value = 1.04157557230905
$ws.Range("F12").Value = 1.049034789847471
$ws.Range("I12").Value = 1.037188443469678
$ws.Range("J12").Value = 1.040324476069615
$ws.Range("K12").Value = 1.041195100609278
$ws.Range("L12").Value = 1.045166826760283
$ws.Range("M12").Value = 1.052598568719622
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033708200767651
$ws.Range("D13").Value = 1.037624238121168
$ws.Range("E13").Value = 1.041618934383063
$ws.Range("F13").Value = 1.049086592702485
$ws.Range("I13").Value = 1.037199359834959
$ws.Range("J13").Value = 1.040352784259686
$ws.Range("K13").Value = 1.041220601598952
$ws.Range("L13").Value = 1.0452004475708
$ws.Range("M13").Value = 1.052640680113123
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033861944896183
$ws.Range("D14").Value = 1.037739769286231
$ws.Range("E14").Value = 1.041761043626646
$ws.Range("F14").Value = 1.049256369949453
$ws.Range("I14").Value = 1.037235098381618
$ws.Range("J14").Value = 1.040445535452341
$ws.Range("K14").Value = 1.041304148849763
$ws.Range("L14").Value = 1.045310617843163
$ws.Range("M14").Value = 1.052778682748915
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033956685492684
$ws.Range("D15").Value = 1.037810962168662
$ws.Range("E15").Value = 1.041848624900794
$ws.Range("F15").Value = 1.049361006934063
$ws.Range("I15").Value = 1.037257095310113
$ws.Range("J15").Value = 1.040502680442082
$ws.Range("K15").Value = 1.041355618439256
$ws.Range("L15").Value = 1.045378504396796
$ws.Range("M15").Value = 1.052863727181154
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03450832592053
$ws.Range("D16").Value = 1.038225491887016
$ws.Range("E16").Value = 1.04235873699253
$ws.Range("F16").Value = 1.04997052103999
$ws.Range("I16").Value = 1.037384778823947
$ws.Range("J16").Value = 1.04083525638118
$ws.Range("K16").Value = 1.041655091638992
$ws.Range("L16").Value = 1.045773739634916
$ws.Range("M16").Value = 1.053358972730338
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034854538613015
$ws.Range("D17").Value = 1.038485651878726
$ws.Range("E17").Value = 1.042679024912835
$ws.Range("F17").Value = 1.050353276081409
$ws.Range("I17").Value = 1.037464564230409
$ws.Range("J17").Value = 1.041043842884809
$ws.Range("K17").Value = 1.041842852477415
$ws.Range("L17").Value = 1.046021752761304
$ws.Range("M17").Value = 1.053669846429196
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035056544242249
$ws.Range("D18").Value = 1.038637447471557
$ws.Range("E18").Value = 1.042865953968741
$ws.Range("F18").Value = 1.050576682550145
$ws.Range("I18").Value = 1.037510990847459
$ws.Range("J18").Value = 1.041165496847957
$ws.Range("K18").Value = 1.041952337123256
$ws.Range("L18").Value = 1.046166447566498
$ws.Range("M18").Value = 1.053851252373479
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035125434033692
$ws.Range("D19").Value = 1.03868921405376
$ws.Range("E19").Value = 1.042929710646993
$ws.Range("F19").Value = 1.050652884105598
$ws.Range("I19").Value = 1.037526802302213
$ws.Range("J19").Value = 1.041206975804127
$ws.Range("K19").Value = 1.041989662914463
$ws.Range("L19").Value = 1.04621579030669
$ws.Range("M19").Value = 1.053913120398721
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03481738647072
$ws.Range("D20").Value = 1.038457734128136
$ws.Range("E20").Value = 1.042644649606258
$ws.Range("F20").Value = 1.05031219436199
$ws.Range("I20").Value = 1.037456015473939
$ws.Range("J20").Value = 1.041021464658537
$ws.Range("K20").Value = 1.04182271092946
$ws.Range("L20").Value = 1.045995139883581
$ws.Range("M20").Value = 1.053636484486099
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033816668212492
$ws.Range("D21").Value = 1.037705746085228
$ws.Range("E21").Value = 1.041719191185709
$ws.Range("F21").Value = 1.049206368191966
$ws.Range("I21").Value = 1.037224579012767
$ws.Range("J21").Value = 1.040418222969479
$ws.Range("K21").Value = 1.041279547645772
$ws.Range("L21").Value = 1.045278173980755
$ws.Range("M21").Value = 1.05273804098215
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033188076892931
$ws.Range("D22").Value = 1.037233390621379
$ws.Range("E22").Value = 1.041138325173444
$ws.Range("F22").Value = 1.048512470718221
$ws.Range("I22").Value = 1.037078071066837
$ws.Range("J22").Value = 1.040038849224012
$ws.Range("K22").Value = 1.04093774898894
$ws.Range("L22").Value = 1.04482769426482
$ws.Range("M22").Value = 1.052173871189808
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033521247748357
$ws.Range("D23").Value = 1.037483752081896
$ws.Range("E23").Value = 1.041446157321925
$ws.Range("F23").Value = 1.048880187638313
$ws.Range("I23").Value = 1.037155831784493
$ws.Range("J23").Value = 1.040239970758251
$ws.Range("K23").Value = 1.041118970075021
$ws.Range("L23").Value = 1.045066472747763
$ws.Range("M23").Value = 1.052472879822404
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034834173721988
$ws.Range("D24").Value = 1.038470348807724
$ws.Range("E24").Value = 1.042660181986289
$ws.Range("F24").Value = 1.050330756953386
$ws.Range("I24").Value = 1.037459878631615
$ws.Range("J24").Value = 1.041031576449883
$ws.Range("K24").Value = 1.041831812130317
$ws.Range("L24").Value = 1.046007164995827
$ws.Range("M24").Value = 1.053651559069499
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036360320897137
$ws.Range("D25").Value = 1.039617142982559
$ws.Range("E25").Value = 1.044073320333562
$ws.Range("F25").Value = 1.052020016198537
$ws.Range("I25").Value = 1.03780835708761
$ws.Range("J25").Value = 1.041949754742332
$ws.Range("K25").Value = 1.042657724194693
$ws.Range("L25").Value = 1.047100073920644
$ws.Range("M25").Value = 1.056118261612939
